# Apply 'storage model issues sorted' edits
#
# summary/zone/generator sheets: numeric updates driven by the revised
# two-unit storage model (PSH + PSH2) on the storage sheet.
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("summary")
$summary.Cells.Item(2,3).Value = 1100
$summary.Cells.Item(3,2).Value = 2100.000000000001
$summary.Cells.Item(3,3).Value = 1600
$summary.Cells.Item(4,2).Value = 2300.000000000001
$summary.Cells.Item(4,3).Value = 2598.148148148149
$summary.Cells.Item(5,3).Value = 4148.148148148149
$summary.Cells.Item(6,3).Value = 300.0000000000003
$summary.Cells.Item(7,3).Value = 500.0000000000002
$summary.Cells.Item(8,3).Value = 0
$summary.Cells.Item(9,3).Value = 0
$summary.Cells.Item(10,3).Value = 0
$summary.Cells.Item(11,2).Value = 2800.000000000001
$summary.Cells.Item(11,3).Value = 850.0000000000002

$zone = $wb.Worksheets.Item("zone")
$zone.Cells.Item(2,4).Value = 1100
$zone.Cells.Item(4,4).Value = 1600
$zone.Cells.Item(5,3).Value = 2100.000000000001
$zone.Cells.Item(6,4).Value = 2598.148148148149
$zone.Cells.Item(7,3).Value = 2300.000000000001
$zone.Cells.Item(8,4).Value = 4148.148148148149
$zone.Cells.Item(10,4).Value = 300.0000000000003
$zone.Cells.Item(12,4).Value = 500.0000000000002
$zone.Cells.Item(14,4).Value = 0
$zone.Cells.Item(16,4).Value = 0
$zone.Cells.Item(18,4).Value = 0
$zone.Cells.Item(20,4).Value = 850.0000000000002
$zone.Cells.Item(21,3).Value = 2800.000000000001

$generator = $wb.Worksheets.Item("generator")
$generator.Cells.Item(6,6).Value = 2100.000000000001
$generator.Cells.Item(9,6).Value = 2300.000000000001
$generator.Cells.Item(30,6).Value = 2800.000000000001

# storage sheet: split the single PSH unit into PSH + PSH2, rows 2-21
$storage = $wb.Worksheets.Item("storage")
$storage.Cells.Item(2,1).Value = "PSH"
$storage.Cells.Item(2,2).Value = 1
$storage.Cells.Item(2,3).Value = 0
$storage.Cells.Item(2,4).Value = 0
$storage.Cells.Item(2,5).Value = 0
$storage.Cells.Item(3,1).Value = "PSH2"
$storage.Cells.Item(3,2).Value = 1
$storage.Cells.Item(3,3).Value = 0
$storage.Cells.Item(3,4).Value = 0
$storage.Cells.Item(3,5).Value = 0
$storage.Cells.Item(4,1).Value = "PSH"
$storage.Cells.Item(4,2).Value = 2
$storage.Cells.Item(4,3).Value = 500
$storage.Cells.Item(4,4).Value = 0
$storage.Cells.Item(4,5).Value = 450
$storage.Cells.Item(5,1).Value = "PSH2"
$storage.Cells.Item(5,2).Value = 2
$storage.Cells.Item(5,3).Value = 500
$storage.Cells.Item(5,4).Value = 0
$storage.Cells.Item(5,5).Value = 450
$storage.Cells.Item(6,1).Value = "PSH"
$storage.Cells.Item(6,2).Value = 3
$storage.Cells.Item(6,3).Value = 1000.000000000001
$storage.Cells.Item(6,4).Value = 0
$storage.Cells.Item(6,5).Value = 900.0000000000007
$storage.Cells.Item(7,1).Value = "PSH2"
$storage.Cells.Item(7,2).Value = 3
$storage.Cells.Item(7,3).Value = 898.1481481481477
$storage.Cells.Item(7,4).Value = 0
$storage.Cells.Item(7,5).Value = 808.3333333333328
$storage.Cells.Item(8,1).Value = "PSH"
$storage.Cells.Item(8,2).Value = 4
$storage.Cells.Item(8,3).Value = 1500.000000000001
$storage.Cells.Item(8,4).Value = 0
$storage.Cells.Item(8,5).Value = 1350.000000000001
$storage.Cells.Item(9,1).Value = "PSH2"
$storage.Cells.Item(9,2).Value = 4
$storage.Cells.Item(9,3).Value = 1398.148148148148
$storage.Cells.Item(9,4).Value = 0
$storage.Cells.Item(9,5).Value = 1258.333333333333
$storage.Cells.Item(10,1).Value = "PSH"
$storage.Cells.Item(10,2).Value = 5
$storage.Cells.Item(10,3).Value = 0
$storage.Cells.Item(10,4).Value = 500
$storage.Cells.Item(10,5).Value = -555.5555555555555
$storage.Cells.Item(11,1).Value = "PSH2"
$storage.Cells.Item(11,2).Value = 5
$storage.Cells.Item(11,3).Value = 0
$storage.Cells.Item(11,4).Value = 500
$storage.Cells.Item(11,5).Value = -555.5555555555555
$storage.Cells.Item(12,1).Value = "PSH"
$storage.Cells.Item(12,2).Value = 6
$storage.Cells.Item(12,3).Value = 0
$storage.Cells.Item(12,4).Value = 1000
$storage.Cells.Item(12,5).Value = -1111.111111111111
$storage.Cells.Item(13,1).Value = "PSH2"
$storage.Cells.Item(13,2).Value = 6
$storage.Cells.Item(13,3).Value = 500
$storage.Cells.Item(13,4).Value = 0
$storage.Cells.Item(13,5).Value = 450
$storage.Cells.Item(14,1).Value = "PSH"
$storage.Cells.Item(14,2).Value = 7
$storage.Cells.Item(14,3).Value = 0
$storage.Cells.Item(14,4).Value = 740.0000000000005
$storage.Cells.Item(14,5).Value = -822.2222222222229
$storage.Cells.Item(15,1).Value = "PSH2"
$storage.Cells.Item(15,2).Value = 7
$storage.Cells.Item(15,3).Value = 0
$storage.Cells.Item(15,4).Value = 359.9999999999995
$storage.Cells.Item(15,5).Value = -399.9999999999995
$storage.Cells.Item(16,1).Value = "PSH"
$storage.Cells.Item(16,2).Value = 8
$storage.Cells.Item(16,3).Value = 0
$storage.Cells.Item(16,4).Value = 190.0000000000007
$storage.Cells.Item(16,5).Value = -211.1111111111119
$storage.Cells.Item(17,1).Value = "PSH2"
$storage.Cells.Item(17,2).Value = 8
$storage.Cells.Item(17,3).Value = 0
$storage.Cells.Item(17,4).Value = 859.9999999999994
$storage.Cells.Item(17,5).Value = -955.555555555555
$storage.Cells.Item(18,1).Value = "PSH"
$storage.Cells.Item(18,2).Value = 9
$storage.Cells.Item(18,3).Value = 0
$storage.Cells.Item(18,4).Value = 0
$storage.Cells.Item(18,5).Value = 0
$storage.Cells.Item(19,1).Value = "PSH2"
$storage.Cells.Item(19,2).Value = 9
$storage.Cells.Item(19,3).Value = 0
$storage.Cells.Item(19,4).Value = 950.0000000000002
$storage.Cells.Item(19,5).Value = -1055.555555555556
$storage.Cells.Item(20,1).Value = "PSH"
$storage.Cells.Item(20,2).Value = 10
$storage.Cells.Item(20,3).Value = 0
$storage.Cells.Item(20,4).Value = 0
$storage.Cells.Item(20,5).Value = 0
$storage.Cells.Item(21,1).Value = "PSH2"
$storage.Cells.Item(21,2).Value = 10
$storage.Cells.Item(21,3).Value = 0
$storage.Cells.Item(21,4).Value = 0
$storage.Cells.Item(21,5).Value = 0
